# fix validate profile, route, and regional
#
# The import-alumni template gains four new "regional" address columns
# (negara, provinsi, kabupaten, kecamatan) inserted right before the
# existing "alamat" column, so the downstream validation/route can map
# Indonesian administrative-region fields independently of the free-text
# address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit was made with K6 selected (this is what ends up recorded as the
# sheet's active cell/selection once the new columns are in place).
$ws.Range("K6").Select() | Out-Null

# Insert 4 blank columns at H:K - this shifts the old H:R ("alamat" .. the
# end) right by four columns to L:V, matching the rest of the sheet.
$ws.Columns("H:K").Insert()

# New header row (row 1) labels for the inserted columns.
$ws.Range("H1").Value = "negara"
$ws.Range("I1").Value = "provinsi"
$ws.Range("J1").Value = "kabupaten"
$ws.Range("K1").Value = "kecamatan"

# Give the new columns sensible widths (data rows 2:4 are left blank, same
# as the other sample/template columns).
$ws.Columns("H").ColumnWidth = 8
$ws.Columns("I:K").ColumnWidth = 11.17
